# Update the date line in the title paragraph.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-05-08 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-05-09 Thursday", 2)

# Update the answer table. The table has 20 rows; every 4th row (1, 5, 9,
# 13, 17) holds the 5 answer cells for that block. Setting Cell.Range.Text
# replaces only the run's text while keeping the existing run/paragraph
# formatting intact.
$table = $d.Tables.Item(1)

$answers = @{
    1  = @("961÷9=106, 7", "608÷3=202, 2", "543÷8=67, 7", "604÷9=67, 1", "993÷3=331, 0")
    5  = @("982÷6=163, 4", "399÷6=66, 3", "858÷5=171, 3", "274÷9=30, 4", "928÷9=103, 1")
    9  = @("428÷3=142, 2", "126÷6=21, 0", "656÷8=82, 0", "634÷9=70, 4", "838÷7=119, 5")
    13 = @("942÷3=314, 0", "796÷3=265, 1", "892÷9=99, 1", "875÷9=97, 2", "961÷2=480, 1")
    17 = @("581÷9=64, 5", "106÷6=17, 4", "900÷3=300, 0", "747÷6=124, 3", "205÷7=29, 2")
}

foreach ($rowIndex in $answers.Keys) {
    $values = $answers[$rowIndex]
    for ($col = 1; $col -le $values.Length; $col++) {
        $cell = $table.Cell($rowIndex, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}

Write-Host "done"
